$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220.8
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = ""
$ws.Range("H11").Value = 107.1875
$ws.Range("I11").Value = 107.1875
$ws.Range("K11").Value = 107.1875
$ws.Range("M11").Value = 32.8125
$ws.Range("H17").Value = 1919.2307
$ws.Range("J17").Value = 1919.2307
$ws.Range("L17").Value = 5757.6921
$ws.Range("N17").Value = -6093.6921
$ws.Range("H33").Value = 511.9091
$ws.Range("I33").Value = 520.5714
$ws.Range("J33").Value = 496.75
$ws.Range("K33").Value = 520.5714
$ws.Range("L33").Value = 496.75
$ws.Range("M33").Value = -291.5714
$ws.Range("N33").Value = -954.75
$ws.Range("H40").Value = 3619.5
$ws.Range("I40").Value = 2528.2666
$ws.Range("K40").Value = 2528.2666
$ws.Range("M40").Value = -2353.2666
$ws.Range("H43").Value = 1724.8889
$ws.Range("I43").Value = 1200
$ws.Range("K43").Value = 1200
$ws.Range("M43").Value = -1131
$ws.Range("H51").Value = 4843.125
$ws.Range("I51").Value = 4000
$ws.Range("J51").Value = 5124.1665
$ws.Range("K51").Value = 4000
$ws.Range("L51").Value = 5124.1665
$ws.Range("N51").Value = -6092.1665
$ws.Range("M51").Value = -3516
$ws.Range("H80").Value = 2494.804
$ws.Range("I80").Value = 1493.28
$ws.Range("J80").Value = 3457.8076
$ws.Range("K80").Value = 4479.84
$ws.Range("L80").Value = 10373.4228
$ws.Range("M80").Value = -3481.84
$ws.Range("N80").Value = -12369.4228
$ws.Range("H83").Value = 2494.804
$ws.Range("I83").Value = 1493.28
$ws.Range("J83").Value = 3457.8076
$ws.Range("K83").Value = 13439.52
$ws.Range("L83").Value = 31120.2684
$ws.Range("M83").Value = -8447.52
$ws.Range("N83").Value = -41104.2684
$ws.Range("H86").Value = 4788634
$ws.Range("I86").Value = 4439.8
$ws.Range("J86").Value = 8775463
$ws.Range("K86").Value = 4439.8
$ws.Range("L86").Value = 8775463
$ws.Range("M86").Value = -3316.8
$ws.Range("N86").Value = -8777709
$ws.Range("H89").Value = 4788634
$ws.Range("I89").Value = 4439.8
$ws.Range("J89").Value = 8775463
$ws.Range("K89").Value = 22199
$ws.Range("L89").Value = 43877315
$ws.Range("M89").Value = -16583
$ws.Range("N89").Value = -43888547
$ws.Range("H99").Value = 674
$ws.Range("J99").Value = 1209
$ws.Range("L99").Value = 3627
$ws.Range("N99").Value = -6623
$ws.Range("H107").Value = 194.55556
$ws.Range("I107").Value = 133.46666
$ws.Range("K107").Value = 133.46666
$ws.Range("M107").Value = 1786.53334
$ws.Range("H111").Value = 3145.2942
$ws.Range("I111").Value = 2183.5454
$ws.Range("J111").Value = 4908.5
$ws.Range("K111").Value = 6550.6362
$ws.Range("L111").Value = 14725.5
$ws.Range("M111").Value = -3483.6362
$ws.Range("N111").Value = -20859.5
$ws.Range("H115").Value = 1314.1666
$ws.Range("I115").Value = 1314.1666
$ws.Range("K115").Value = 3942.4998
$ws.Range("M115").Value = -2375.4998
$ws.Range("H116").Value = 3666.6667
$ws.Range("I116").Value = 3000
$ws.Range("K116").Value = 3000
$ws.Range("M116").Value = 442
$ws.Range("H132").Value = 37498.285
$ws.Range("I132").Value = 1787.7
$ws.Range("J132").Value = 126774.75
$ws.Range("K132").Value = 5363.1
$ws.Range("L132").Value = 380324.25
$ws.Range("M132").Value = -2833.1
$ws.Range("N132").Value = -385384.25
$ws.Range("H137").Value = 2922.1333
$ws.Range("I137").Value = 1793.25
$ws.Range("J137").Value = 3332.6365
$ws.Range("K137").Value = 5379.75
$ws.Range("L137").Value = 9997.9095
$ws.Range("M137").Value = -2829.75
$ws.Range("N137").Value = -15097.9095
$ws.Range("H138").Value = 5599.6
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 5599.6
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 16798.8
$ws.Range("M138").Value = ""
$ws.Range("N138").Value = -27078.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3492.25
$ws.Range("I45").Value = 1984.8
$ws.Range("J45").Value = 6004.6665
$ws.Range("K45").Value = 1984.8
$ws.Range("L45").Value = 6004.6665
$ws.Range("M45").Value = -1607.8
$ws.Range("N45").Value = -6758.6665
$ws.Range("H61").Value = 6210.3335
$ws.Range("I61").Value = 4932
$ws.Range("K61").Value = 4932
$ws.Range("M61").Value = -4720
$ws.Range("H97").Value = 4274225.5
$ws.Range("I97").Value = 703.0952
$ws.Range("K97").Value = 703.0952
$ws.Range("M97").Value = -207.0952
$ws.Range("H102").Value = 20834894
$ws.Range("I102").Value = 1640.3572
$ws.Range("K102").Value = 1640.3572
$ws.Range("M102").Value = -18.35719999999992
$ws.Range("H126").Value = 7200
$ws.Range("I126").Value = 7200
$ws.Range("K126").Value = 21600
$ws.Range("M126").Value = -19130
$ws.Range("H132").Value = 2942.4138
$ws.Range("I132").Value = 2015.3636
$ws.Range("J132").Value = 5856
$ws.Range("K132").Value = 6046.0908
$ws.Range("L132").Value = 17568
$ws.Range("M132").Value = -3516.0908
$ws.Range("N132").Value = -22628
$ws.Range("H136").Value = 6210.3335
$ws.Range("I136").Value = 4932
$ws.Range("K136").Value = 14796
$ws.Range("M136").Value = -12246

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3004.7058
$ws.Range("I99").Value = 2823.8462
$ws.Range("K99").Value = 2823.8462
$ws.Range("M99").Value = -1325.8462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3101.7673
$ws.Range("I31").Value = 2111.5925
$ws.Range("J31").Value = 4772.6875
$ws.Range("K31").Value = 2111.5925
$ws.Range("L31").Value = 4772.6875
$ws.Range("M31").Value = -1816.5925
$ws.Range("N31").Value = -5362.6875
$ws.Range("H34").Value = 3101.7673
$ws.Range("I34").Value = 2111.5925
$ws.Range("J34").Value = 4772.6875
$ws.Range("K34").Value = 2111.5925
$ws.Range("L34").Value = 4772.6875
$ws.Range("M34").Value = -1909.5925
$ws.Range("N34").Value = -5176.6875
$ws.Range("H36").Value = 999.5
$ws.Range("I36").Value = 999.5
$ws.Range("K36").Value = 999.5
$ws.Range("M36").Value = -611.5
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").Value = ""
$ws.Range("H40").Value = 999.5
$ws.Range("I40").Value = 999.5
$ws.Range("K40").Value = 999.5
$ws.Range("M40").Value = -839.5
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = ""
$ws.Range("H99").Value = 11817734
$ws.Range("I99").Value = 4071373
$ws.Range("K99").Value = 4071373
$ws.Range("M99").Value = -4069875
$ws.Range("H105").Value = 13163539
$ws.Range("I105").Value = 1283.3334
$ws.Range("J105").Value = 62522000
$ws.Range("K105").Value = 1283.3334
$ws.Range("L105").Value = 62522000
$ws.Range("M105").Value = 463.6666
$ws.Range("N105").Value = -62525494
$ws.Range("H122").Value = 303426.5
$ws.Range("J122").Value = 5254.9443
$ws.Range("L122").Value = 15764.8329
$ws.Range("N122").Value = -20664.8329
$ws.Range("H126").Value = 11817734
$ws.Range("I126").Value = 4071373
$ws.Range("K126").Value = 12214119
$ws.Range("M126").Value = -12211649
$ws.Range("H141").Value = 107230.695
$ws.Range("J141").Value = 107230.695
$ws.Range("L141").Value = 107230.695
$ws.Range("N141").Value = -117590.695

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 55651.293
$ws.Range("J70").Value = 11355.857
$ws.Range("L70").Value = 11355.857
$ws.Range("N70").Value = -11895.857
$ws.Range("H73").Value = 55651.293
$ws.Range("J73").Value = 11355.857
$ws.Range("L73").Value = 11355.857
$ws.Range("N73").Value = -13227.857
$ws.Range("H132").Value = 2698.375
$ws.Range("I132").Value = 1841.3334
$ws.Range("K132").Value = 5524.0002
$ws.Range("M132").Value = -2994.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 401.75
$ws.Range("I16").Value = 264
$ws.Range("J16").Value = 677.25
$ws.Range("K16").Value = 264
$ws.Range("L16").Value = 677.25
$ws.Range("M16").Value = -94
$ws.Range("N16").Value = -1017.25
$ws.Range("H22").Value = 1162.6
$ws.Range("I22").Value = 975.2857
$ws.Range("K22").Value = 975.2857
$ws.Range("M22").Value = -680.2857
$ws.Range("H27").Value = 1162.6
$ws.Range("I27").Value = 975.2857
$ws.Range("K27").Value = 975.2857
$ws.Range("M27").Value = -868.2857
$ws.Range("H36").Value = 57749
$ws.Range("J36").Value = 57749
$ws.Range("L36").Value = 57749
$ws.Range("N36").Value = -58873
$ws.Range("H93").Value = 1863.3667
$ws.Range("J93").Value = 2171.2727
$ws.Range("L93").Value = 2171.2727
$ws.Range("N93").Value = -4667.2727

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 37039404
$ws.Range("J62").Value = 37039404
$ws.Range("L62").Value = 37039404
$ws.Range("N62").Value = -37040652
$ws.Range("H64").Value = 49492.2
$ws.Range("I64").Value = 48500
$ws.Range("J64").Value = 50153.668
$ws.Range("K64").Value = 48500
$ws.Range("L64").Value = 50153.668
$ws.Range("M64").Value = -48252
$ws.Range("N64").Value = -50649.668
$ws.Range("H65").Value = 37039404
$ws.Range("J65").Value = 37039404
$ws.Range("L65").Value = 185197020
$ws.Range("N65").Value = -185203260
$ws.Range("H67").Value = 49492.2
$ws.Range("I67").Value = 48500
$ws.Range("J67").Value = 50153.668
$ws.Range("K67").Value = 48500
$ws.Range("L67").Value = 50153.668
$ws.Range("M67").Value = -47642
$ws.Range("N67").Value = -51869.668
$ws.Range("H122").Value = 375058.28
$ws.Range("I122").Value = 559533.6
$ws.Range("K122").Value = 559533.6
$ws.Range("M122").Value = -1676150.8
